$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet "My Data" -> "People"
$ws.Name = "People"

# Fix header typo: "Birthdate" -> "BirthDate"
$ws.Range("B1").Value = "BirthDate"

# Update existing birthdate values to include time component
$ws.Range("B2").Value = "18/12/2002 00:00:00"
$ws.Range("B3").Value = "19/03/2000 00:00:00"

# Add a new birthdate value for row 4 (Jonathan)
$ws.Range("B4").Value = "20/05/1999 00:00:00"
